# Finished with temporary admin page
#
# The "filename" column is renamed to "image", and the previously-empty
# image cell for the "Lemon Drip" row (row 6) is filled in with its
# filename. Order matters here: writing H6 before H1 reproduces the same
# shared-string table ordering as the source edit (lemon_drip.jpg ends up
# right before "image" at the tail of the shared strings table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Lemon Drip") previously had no filename/image value - give it one.
$ws.Range("H6").Value = "lemon_drip.jpg"

# Rename the "filename" header to "image".
$ws.Range("H1").Value = "image"

# Match the author's final on-screen selection/scroll state.
[void]$ws.Range("H1").Select()
